# Auto-generated edit script: updates market-price derived columns (H:N)
# across multiple sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 2112.75
$ws.Range("I80").Value = 444
$ws.Range("J80").Value = 2351.1428
$ws.Range("K80").Value = 1332
$ws.Range("L80").Value = 7053.428400000001
$ws.Range("M80").Value = -334
$ws.Range("N80").Value = -9049.428400000001
# Row 83
$ws.Range("H83").Value = 2112.75
$ws.Range("I83").Value = 444
$ws.Range("J83").Value = 2351.1428
$ws.Range("K83").Value = 3996
$ws.Range("L83").Value = 21160.2852
$ws.Range("M83").Value = 996
$ws.Range("N83").Value = -31144.2852
# Row 129
$ws.Range("H129").Value = 1634
$ws.Range("I129").Value = 1249.5
$ws.Range("J129").Value = 1710.9
$ws.Range("K129").Value = 3748.5
$ws.Range("L129").Value = 5132.700000000001
$ws.Range("M129").Value = 1251.5
$ws.Range("N129").Value = -15132.7
# Row 132
$ws.Range("H132").Value = 1283.7059
$ws.Range("I132").Value = 935.8461
$ws.Range("J132").Value = 2414.25
$ws.Range("K132").Value = 2807.5383
$ws.Range("L132").Value = 7242.75
$ws.Range("M132").Value = -277.5383000000002
$ws.Range("N132").Value = -12302.75
# Row 137
$ws.Range("H137").Value = 2000
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -3450
$ws.Range("N137").Value = -11100

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 12
# Row 28
$ws.Range("H28").Value = 8777.5
$ws.Range("I28").Value = 8777.5
$ws.Range("K28").Value = 8777.5
$ws.Range("M28").Value = -8585.5
# Row 61
$ws.Range("H61").Value = 3250
$ws.Range("I61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("M61").Value = -1288
# Row 99
$ws.Range("H99").Value = 8777.5
$ws.Range("I99").Value = 8777.5
$ws.Range("K99").Value = 8777.5
$ws.Range("M99").Value = -5782.5
# Row 122
$ws.Range("H122").Value = 3242
$ws.Range("I122").Value = 2582.75
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7748.25
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -5298.25
$ws.Range("N122").Value = -19900
# Row 132
$ws.Range("H132").Value = 2030.1613
$ws.Range("I132").Value = 1025.5714
$ws.Range("J132").Value = 4139.8
$ws.Range("K132").Value = 3076.7142
$ws.Range("L132").Value = 12419.4
$ws.Range("M132").Value = -546.7142000000003
$ws.Range("N132").Value = -17479.4
# Row 136
$ws.Range("H136").Value = 3250
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 15
# Row 86
$ws.Range("H86").Value = 3188.2856
$ws.Range("J86").Value = 3660.5
$ws.Range("L86").Value = 3660.5
$ws.Range("N86").Value = -5906.5
# Row 89
$ws.Range("H89").Value = 3188.2856
$ws.Range("J89").Value = 3660.5
$ws.Range("L89").Value = 18302.5
$ws.Range("N89").Value = -29534.5
# Row 94
$ws.Range("H94").Value = 531.8182
$ws.Range("I94").Value = 516.6667
$ws.Range("K94").Value = 516.6667
$ws.Range("M94").Value = -65.66669999999999
# Row 134
$ws.Range("H134").Value = 3055.5557
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4014.7778
$ws.Range("I31").Value = 3095.875
$ws.Range("K31").Value = 3095.875
$ws.Range("M31").Value = -2800.875
# Row 34
$ws.Range("H34").Value = 4014.7778
$ws.Range("I34").Value = 3095.875
$ws.Range("K34").Value = 3095.875
$ws.Range("M34").Value = -2893.875
# Row 58
$ws.Range("H58").Value = 2209.75
$ws.Range("I58").Value = 2209.75
$ws.Range("K58").Value = 2209.75
$ws.Range("M58").Value = -2006.75
# Row 74
$ws.Range("H74").Value = 53313.332
$ws.Range("J74").Value = 53313.332
$ws.Range("L74").Value = 53313.332
$ws.Range("N74").Value = -55061.332
# Row 77
$ws.Range("H77").Value = 53313.332
$ws.Range("J77").Value = 53313.332
$ws.Range("L77").Value = 159939.996
$ws.Range("N77").Value = -168675.996
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
# Row 132
$ws.Range("H132").Value = 2682.4614
$ws.Range("I132").Value = 2287.4
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 6862.200000000001
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -4332.200000000001
$ws.Range("N132").Value = -17057.9999
# Row 136
$ws.Range("H136").Value = 2209.75
$ws.Range("I136").Value = 2209.75
$ws.Range("K136").Value = 6629.25
$ws.Range("M136").Value = -4079.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3013.5334
$ws.Range("I132").Value = 2517.1667
$ws.Range("K132").Value = 7551.500100000001
$ws.Range("M132").Value = -5021.500100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 1948.75
$ws.Range("I40").Value = 1948.75
$ws.Range("K40").Value = 1948.75
$ws.Range("M40").Value = -1812.75
# Row 46
$ws.Range("H46").Value = 2691.923
$ws.Range("I46").Value = 1999.5
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1999.5
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1811.5
$ws.Range("N46").Value = -5376
# Row 68
$ws.Range("H68").Value = 1802.125
$ws.Range("I68").Value = 1083.6
$ws.Range("J68").Value = 2999.6667
$ws.Range("K68").Value = 1083.6
$ws.Range("L68").Value = 2999.6667
$ws.Range("M68").Value = -334.5999999999999
$ws.Range("N68").Value = -4497.6667
# Row 71
$ws.Range("H71").Value = 1802.125
$ws.Range("I71").Value = 1083.6
$ws.Range("J71").Value = 2999.6667
$ws.Range("K71").Value = 5418
$ws.Range("L71").Value = 14998.3335
$ws.Range("M71").Value = -1674
$ws.Range("N71").Value = -22486.3335
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
# Row 122
$ws.Range("H122").Value = 2929.8
$ws.Range("I122").Value = 2929.8
$ws.Range("K122").Value = 8789.400000000001
$ws.Range("M122").Value = -6339.400000000001
# Row 132
$ws.Range("H132").Value = 2123.5
$ws.Range("I132").Value = 1157.8
$ws.Range("K132").Value = 3473.4
$ws.Range("M132").Value = -943.3999999999996

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2749.5
$ws.Range("I126").Value = 2749.5
$ws.Range("K126").Value = 8248.5
$ws.Range("M126").Value = -5778.5
# Row 132
$ws.Range("H132").Value = 2451.9678
$ws.Range("I132").Value = 926.3125
$ws.Range("J132").Value = 4079.3333
$ws.Range("K132").Value = 2778.9375
$ws.Range("L132").Value = 12237.9999
$ws.Range("M132").Value = -248.9375
$ws.Range("N132").Value = -17297.9999
